$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (A: 15.42578125 -> 16.42578125, B: 14.7109375 -> 15.7109375)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 14.833333333333334

# New values for columns A and B (rows 1-32)
$aValues = @(-0.20466831670470498,-0.18243722887989744,-0.13221075427419393,-0.1238999308102624,-0.1202832684299624,-0.017595135142698837,-0.0074528080214046888,0.011695910230349327,0.013703936640774295,0.01571946874517316,-0.0079458345732694724,-0.0044456672913328887,-0.00094415851635343273,-0.0090786694607558616,-0.0080510661706076903,-0.0060332435218155567,-0.0040033210823118281,-0.016103890361190309,-0.012091285579486577,-0.0080165752770273713,-0.0040056327046933404,-0.045703642177837622,-0.040492872234080046,-0.020097790987785302,-0.15486813418189271,-0.094605589169219684,-0.091946605108276902,-0.089004316358983537,-0.081359199151328276,-0.021172204833151742,-0.014023476368576482,-0.0040012161340090557)
$bValues = @(0.20454111920705742,0.18191439642029472,0.13189993075244821,0.12328326839501003,0.11816579753364032,0.01745280793681836,0.0074315687641841599,-0.01170393668149794,-0.013719468784762157,-0.015720402735226457,0.007945667243447474,0.0044441584678138746,0.00094287572773321671,0.0090510661350284849,0.0080332434802166119,0.0060033210398859893,0.0039999999464530589,0.016091285556210977,0.012016575252087769,0.0080056326794704091,0.0039999999745612413,0.045492872195074696,0.040097790863971916,0.019999999874558583,0.15454399639217975,0.094446605059534505,0.091004316302895294,0.088359199060381144,0.081172204449636354,0.021023476271128771,0.01400121601981752,0.0039999999183280011)

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}
